$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Language" column (D) first, then the "Title" column (B),
# so that "ContentType" (currently column C) shifts left into column B,
# leaving just Path (A) and ContentType (B).
$ws.Columns.Item(4).Delete()
$ws.Columns.Item(2).Delete()

# Update the selection to match the recorded state (A6)
$ws.Range("A6").Select()

